$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.036961
$ws.Range("H2").Value = 0.110883
$ws.Range("I2").Value = 0.1786608532565087
$ws.Range("J2").Value = 0.1786608532565087
$ws.Range("M2").Value = 0.01451666666666667
$ws.Range("N2").Value = 0.04355
$ws.Range("O2").Value = 0.09504707612932513
$ws.Range("P2").Value = 0.09504707612932513
$ws.Range("Q2").Value = 0.0005365505166666666
$ws.Range("R2").Value = 0.00482895465
$ws.Range("S2").Value = 0.01698119172080156
$ws.Range("T2").Value = 0.01698119172080156
$ws.Range("G3").Value = 0.036961
$ws.Range("H3").Value = 0.110883
$ws.Range("I3").Value = 0.1786608532565087
$ws.Range("J3").Value = 0.1786608532565087
$ws.Range("O3").Value = 0.5165388459909994
$ws.Range("P3").Value = 0.5165388459909994
$ws.Range("Q3").Value = 0.002915914891666667
$ws.Range("R3").Value = 0.026243234025
$ws.Range("S3").Value = 0.09228527096488427
$ws.Range("T3").Value = 0.09228527096488427
$ws.Range("G4").Value = 0.036961
$ws.Range("H4").Value = 0.110883
$ws.Range("I4").Value = 0.1786608532565087
$ws.Range("J4").Value = 0.1786608532565087
$ws.Range("O4").Value = 0.3884140778796754
$ws.Range("P4").Value = 0.3884140778796754
$ws.Range("Q4").Value = 0.002192637403
$ws.Range("R4").Value = 0.019733736627
$ws.Range("S4").Value = 0.06939439057082282
$ws.Range("T4").Value = 0.06939439057082282
$ws.Range("H5").Value = 0.384071
$ws.Range("I5").Value = 0.6188365445657183
$ws.Range("J5").Value = 0.6188365445657182
$ws.Range("M5").Value = 0.01451666666666667
$ws.Range("N5").Value = 0.04355
$ws.Range("O5").Value = 0.09504707612932513
$ws.Range("P5").Value = 0.09504707612932513
$ws.Range("Q5").Value = 0.001858476894444444
$ws.Range("R5").Value = 0.01672629205
$ws.Range("S5").Value = 0.05881860416294633
$ws.Range("T5").Value = 0.05881860416294633
$ws.Range("H6").Value = 0.384071
$ws.Range("I6").Value = 0.6188365445657183
$ws.Range("J6").Value = 0.6188365445657182
$ws.Range("O6").Value = 0.5165388459909994
$ws.Range("P6").Value = 0.5165388459909994
$ws.Range("R6").Value = 0.090900003925
$ws.Range("S6").Value = 0.3196531145870338
$ws.Range("T6").Value = 0.3196531145870338
$ws.Range("H7").Value = 0.384071
$ws.Range("I7").Value = 0.6188365445657183
$ws.Range("J7").Value = 0.6188365445657182
$ws.Range("O7").Value = 0.3884140778796754
$ws.Range("P7").Value = 0.3884140778796754
$ws.Range("R7").Value = 0.06835273179899999
$ws.Range("S7").Value = 0.2403648258157381
$ws.Range("T7").Value = 0.2403648258157381
$ws.Range("I8").Value = 0.202502602177773
$ws.Range("J8").Value = 0.202502602177773
$ws.Range("M8").Value = 0.01451666666666667
$ws.Range("N8").Value = 0.04355
$ws.Range("O8").Value = 0.09504707612932513
$ws.Range("P8").Value = 0.09504707612932513
$ws.Range("Q8").Value = 0.0006081515555555555
$ws.Range("R8").Value = 0.005473363999999999
$ws.Range("S8").Value = 0.01924728024557723
$ws.Range("T8").Value = 0.01924728024557723
$ws.Range("I9").Value = 0.202502602177773
$ws.Range("J9").Value = 0.202502602177773
$ws.Range("O9").Value = 0.5165388459909994
$ws.Range("P9").Value = 0.5165388459909994
$ws.Range("S9").Value = 0.1046004604390813
$ws.Range("T9").Value = 0.1046004604390813
$ws.Range("I10").Value = 0.202502602177773
$ws.Range("J10").Value = 0.202502602177773
$ws.Range("O10").Value = 0.3884140778796754
$ws.Range("P10").Value = 0.3884140778796754
$ws.Range("S10").Value = 0.07865486149311446
$ws.Range("T10").Value = 0.07865486149311446
